# Update the "Rules" sheet: cell E8 held the greeting text shown for rule R10.
# It previously read "Good Morning"; replace it with "GIT UPDATE" and leave
# that cell selected, matching how the change was made interactively in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("E8")
$cell.Value = "GIT UPDATE"
$cell.Select()
